# "Add init message about car for TL" — add an x/y-coordinate "init
# message" pair of columns (x-coord / y-coord) to the Cars sheet so the
# traffic-light (TL) agents have each car's starting position.

$wb = $excel.ActiveWorkbook

$wsTL   = $wb.Worksheets.Item("TrafficLights")
$wsCars = $wb.Worksheets.Item("Cars")

# --- Cars sheet: new G/H columns (x-coord / y-coord) ------------------
$wsCars.Activate()

$wsCars.Range("G1").Value = "x-coord"
$wsCars.Range("H1").Value = "y-coord"

$wsCars.Range("G2").Value = 0
$wsCars.Range("H2").Value = 1010

$wsCars.Range("G3").Value = 10
$wsCars.Range("H3").Value = 2010

$wsCars.Range("G4").Value = 0
$wsCars.Range("H4").Value = 10

$wsCars.Range("G5").Value = 5
$wsCars.Range("H5").Value = 10

# Widen the new x-coord column like Excel does after typing into it
# (target stored width 16.85546875; the interop column-width grid here
# snaps to the nearest 1/6 character, so 16.0 is the closest settable
# value that lands on it).
$wsCars.Columns.Item(7).ColumnWidth = 16.0

# Page setup touched (portrait / A4-ish "9" paper) as part of the same edit.
$wsCars.PageSetup.PaperSize = 9
$wsCars.PageSetup.Orientation = 1

# Final selection left on the sheet after the edit.
$null = $wsCars.Range("H5").Select()

# --- TrafficLights sheet: selection moved to its own new-ish G2:H2 ----
$wsTL.Activate()
$null = $wsTL.Range("G2:H2").Select()

# Leave the workbook on the Cars tab (matches tabSelected/activeTab).
$wsCars.Activate()
